$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Samsung Mikrowelle row (previously row 12); every row below it
# shifts up by one to fill the gap.
$ws.Rows(12).Delete()

# The crawl re-ran later the same day; stamp every remaining data row
# (rows 2-33 after the deletion) with the new crawl timestamp.
$ws.Range("O2:O33").Value = "2022-09-17 21:00:03"
